# Simulated Wild Card round and logged it
# Update the "R" (Road) row target-depth stats on both the OFF and DEF
# sheets to reflect the additional Wild Card game.

$wb = $excel.ActiveWorkbook

$off = $wb.Worksheets.Item("OFF")
$off.Range("B3").Value = 216
$off.Range("C3").Value = 154
$off.Range("D3").Value = 41
$off.Range("E3").Value = 16

$def = $wb.Worksheets.Item("DEF")
$def.Range("B3").Value = 240
$def.Range("C3").Value = 166
$def.Range("D3").Value = 67
$def.Range("E3").Value = 30
